# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to the freshly scraped values from the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 255
$wsExhibit.Range("F5").Value = 6533
$wsExhibit.Range("F6").Value = 5307
$wsExhibit.Range("F10").Value = 63
$wsExhibit.Range("F11").Value = 227
$wsExhibit.Range("F12").Value = 38

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 255
$wsAll.Range("F5").Value = 6533
$wsAll.Range("F6").Value = 5307
$wsAll.Range("F10").Value = 63
$wsAll.Range("F11").Value = 227
$wsAll.Range("F14").Value = 38
